# Fix typo: change 'celltypes' to 'cell types' in the KG worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("J2").Value = "cell types"
